$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("线束表转化")
$ws1.Range("B35").Value = "test"
